# "Tried to implement Penality Reward System (unfinished)"
# The author removed the penalty-period data rows (the weeks/month that
# fell inside the March 2024 "penalty" window) from both the Weekly
# Quantity and Monthly Trend sheets. Everything below shifts up to close
# the gap, and the sheet `dimension` shrinks accordingly.

$wb = $excel.ActiveWorkbook

# --- "Weekly Quantity" sheet: drop the two weekly rows for March 2024 ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

# Row 6 (2024-03-17, qty 250) - delete highest row index first so the
# other target row's index doesn't shift before we get to it.
$wsWeekly.Rows.Item(6).Delete()
# Row 5 (2024-03-10, qty 280)
$wsWeekly.Rows.Item(5).Delete()

# --- "Monthly Trend" sheet: drop the March 2024 monthly row ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Row 3 (2024-03-31, qty 530)
$wsMonthly.Rows.Item(3).Delete()
